# Refresh cryptos.xlsx price/volume figures (and two swapped rows) to match the
# latest scrape, mirroring the supplied OOXML diff cell-by-cell.
#
# Price cells (column D) are written with a leading single-quote so Excel keeps
# them as literal text (e.g. "9.00", "0.0230") instead of silently reparsing them
# as numbers and dropping meaningful trailing/leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "62.313.16"
$ws.Range("E2").Value = "  +1.51%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.425.18"
$ws.Range("E3").Value = "  +1.93%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").Value = "'563.78"
$ws.Range("E5").Value = "  +2.12%  "

# Row 6: Solana
$ws.Range("D6").Value = "'144.85"
$ws.Range("E6").Value = "  +3.58%  "

# Row 8: XRP
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  +1.99%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.423.41"
$ws.Range("E9").Value = "  +1.81%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +2.07%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -1.51%  "

# Row 12: Toncoin
$ws.Range("E12").Value = "  +0.86%  "

# Row 13: Cardano
$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "  +0.75%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'26.10"
$ws.Range("E14").Value = "  +1.98%  "

# Row 15: ShibaInu
$ws.Range("D15").Value = "'0.0000179"

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.862.00"
$ws.Range("E16").Value = "  +1.90%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "61.936.33"
$ws.Range("E17").Value = "  +0.98%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.423.68"
$ws.Range("E18").Value = "  +1.81%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'11.34"
$ws.Range("E19").Value = "  +3.09%  "

# Row 20: Polkadot
$ws.Range("E20").Value = "  +1.24%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'325.34"
$ws.Range("E21").Value = "  +1.30%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'6.78"
$ws.Range("E22").Value = "  +0.93%  "

# Row 23: Dai
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'65.66"
$ws.Range("E24").Value = "  +1.89%  "

# Row 25: SuiNetwork
$ws.Range("E25").Value = "  -2.53%  "

# Row 26: Aptos
$ws.Range("D26").Value = "'9.00"
$ws.Range("E26").Value = "  +1.11%  "

# Row 27: Bittensor
$ws.Range("D27").Value = "'589.82"
$ws.Range("E27").Value = "  +14.42%  "

# Row 28: PEPE
$ws.Range("D28").Value = "0.0₃0952"
$ws.Range("E28").Value = "  +5.80%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.530.39"
$ws.Range("E29").Value = "  +1.35%  "

# Row 30: WrappedeETH
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.09%  "

# Row 31: Fetch.AI
$ws.Range("E31").Value = "  +5.87%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = "'8.30"
$ws.Range("E32").Value = "  +1.11%  "

# Row 33: Kaspa
$ws.Range("E33").Value = "  +0.49%  "

# Row 34: PancakeSwap
$ws.Range("E34").Value = "  +2.10%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +1.11%  "

# Row 36: RenderToken
$ws.Range("D36").Value = "'5.77"
$ws.Range("E36").Value = "  +4.80%  "

# Row 37: FirstDigitalUSD
$ws.Range("E37").Value = "  +0.03%  "

# Row 38: NEARProtocol
$ws.Range("D38").Value = "'4.82"
$ws.Range("E38").Value = "  +2.63%  "

# Row 39: Monero
$ws.Range("D39").Value = "'154.21"
$ws.Range("E39").Value = "  +5.10%  "

# Row 40: PolygonEcosystemToken
$ws.Range("E40").Value = "  +1.40%  "

# Row 41: EthereumClassic
$ws.Range("D41").Value = "'18.73"
$ws.Range("E41").Value = "  +1.01%  "

# Row 42: Stacks
$ws.Range("E42").Value = "  -2.31%  "

# Row 43: USDe
$ws.Range("E43").Value = "  -0.16%  "

# Row 44: dogwifhat
$ws.Range("E44").Value = "  +8.97%  "

# Row 45: Aave
$ws.Range("D45").Value = "'150.66"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46: Filecoin
$ws.Range("D46").Value = "'3.66"
$ws.Range("E46").Value = "  +1.48%  "

# Row 47: Hedera
$ws.Range("D47").Value = "'0.0541"
$ws.Range("E47").Value = "  +2.47%  "

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = "'20.49"
$ws.Range("E48").Value = "  +4.24%  "

# Row 49: Mantle
$ws.Range("D49").Value = "'0.594"
$ws.Range("E49").Value = "  +2.26%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  +2.26%  "

# Row 51: VeChain
$ws.Range("D51").Value = "'0.0230"
$ws.Range("E51").Value = "  +1.78%  "
